$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.946.61'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.893.07'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.02%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7739'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.79'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.33%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3138'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.93%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.74'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07356'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08066'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7714'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.503'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.875.02'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.25'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.223'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.87%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.930.56'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.00'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.94%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '246.75'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007862'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.98%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.149'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.144.33'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.27%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1581'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.439'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.99'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.77'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.025'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.423'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.65%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.541'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.468'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05568'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.062'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.41%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.239'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7522'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.86%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.682'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.52%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01931'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.06%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.48'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4472'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.101.51'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.015'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.33%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.889'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.38'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.60%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.827'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.77%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.538'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.37%  '
